$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 25.23
$ws.Range("A3").Value = "a"
$ws.Range("B3").Value = 10.220000000000001
$ws.Range("F6").Select()
